# "updated templates (removed id column)"
#
# Resource_Template.xlsx originally had columns:
#   A=Id, B=Resource Category, C=Name, D=Description, E=Link
#
# Delete the "Id" column entirely so the remaining columns shift left:
#   A=Resource Category, B=Name, C=Description, D=Link
# Excel automatically re-points the header cells/shared-string refs, the
# used-range dimension, the row's column span, and the catch-all
# "don't type outside the template" data-validation range (which always
# starts at the column just past the last real column) when a column is
# deleted this way.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resource_Template")

$ws.Columns("A:A").Delete()

# Post-edit cursor position recorded in the template.
$ws.Range("D6").Select()
